# Fix formatting bugs: correct a typo, reformat an amount as currency text,
# and add a new "Income" row to the budget table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fix description typo and change amount to a currency-formatted string.
$ws.Range("B2").Value = "asfd"
# Leading apostrophe forces Excel to keep this as literal text instead of
# auto-converting "$12.00" into a formatted number.
$ws.Range("C2").Value = "'$12.00"

# Row 3: new income line.
$ws.Range("A3").Value = "Income"
$ws.Range("B3").Value = "efa"
$ws.Range("C3").Value = "'$1.00"
